# Update the interactive map data: the oldest pending case
# ("6303" - BILBAO, FRANCISCO 2362, row 78) has been resolved/removed from
# the report, and a brand-new case ("-517" - Av Dorrego 2721) was appended
# at the bottom. Net effect on the sheet: delete row 78 entirely, which
# shifts every subsequent row (79-89) up by one, turning the former 12-row
# block (rows 78-89) into an 11-row block (rows 78-88) and shrinking the
# used range from A1:P89 to A1:P88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(78).Delete()
